$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- 1. Re-style the existing row 2 (gets a top/bottom separator look) ---
$ws.Range("A2:E2").Borders(8).LineStyle = 1
$ws.Range("A2:E2").Borders(8).Weight = -4138
$ws.Range("A2:E2").Borders(9).LineStyle = 1
$ws.Range("A2:E2").Borders(9).Weight = 2

# --- 2. Insert three new rows below row 2 (rows 3, 4, 5) ---
$ws.Rows("3:5").Insert()

# --- 3. Row 3 content ---
$ws.Range("A3").Value = "SCRIPT/P01P04A/um0601.ssb"
$ws.Range("B3").Value = 382
$ws.Range("C3").Value = " I have really taken to this shop."
$ws.Range("D3").Value = " Мне очень понравилось это место."
$ws.Range("E3").Value = " Íîå ïœåîû ðïîñàâéìïòû üóï íåòóï."

# --- 4. Row 4 content ---
$ws.Range("A4").Value = "SCRIPT/P01P04A/um0705.ssb"
$ws.Range("B4").Value = 382
$ws.Range("C4").Value = " I have really taken to this shop."
$ws.Range("D4").Value = " Мне очень понравилось это место."
$ws.Range("E4").Value = " Íîå ïœåîû ðïîñàâéìïòû üóï íåòóï."

# --- 5. Row 5 content ---
$ws.Range("A5").Value = "SCRIPT/P01P04A/um0712.ssb"
$ws.Range("B5").Value = 363
$ws.Range("C5").Value = " Recycling is a trend of the time."
$ws.Range("D5").Value = " Переработка - это новый тренд."
$ws.Range("E5").Value = " Ðåñåñàáïóëà - üóï îïâúê óñåîä."

# --- 6. Formatting for rows 3 and 4 (thin top + thin bottom borders, wrap text) ---
$ws.Range("A3:E4").WrapText = $true
$ws.Range("A3:E4").Borders(8).LineStyle = 1
$ws.Range("A3:E4").Borders(9).LineStyle = 1

# --- 7. Row heights for the new rows ---
$ws.Rows("3:5").RowHeight = 43.2

# --- 8. Sheet view: scroll + selection like the target state ---
$ws.Range("D5").Select()
$excel.ActiveWindow.ScrollRow = 4
